# Apply updated values to rows 4, 5, 6 for the Pulping_machines optimization sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("H", "K", "L", "P", "Q", "R", "S", "U", "AC", "AD", "AE")
$newValues = @{
    "H"  = 1.074999854899943
    "K"  = 0.0288218287169002
    "L"  = 0.0009009608766064048
    "P"  = 0.1426123455166817
    "Q"  = 0.00008039055683184415
    "R"  = 0.003843874612357467
    "S"  = 0.000307947862893343
    "U"  = 0.03848720947280526
    "AC" = 0.8509804988462096
    "AD" = 1.263248519040644
    "AE" = 0.7998103252612054
}

foreach ($row in 4..6) {
    foreach ($col in $columns) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}

$wb.Save()
